$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.207.80"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.324.91"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +0.61%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "99.67"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +2.19%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.45"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.29%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  -0.71%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "17.82"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "2.686.21"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "2.353.25"
$ws.Range("E16").Value = "  +1.94%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.800"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "43.126.01"
$ws.Range("E18").Value = "  +0.24%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +0.34%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.26"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "238.36"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +4.64%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.47"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "25.45"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "168.31"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "34.58"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("E30").Value = "  +0.17%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.04"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.16%  "
$ws.Range("E32").Value = "  +4.13%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.73"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.92%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "17.61"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "2.000.73"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("E45").Value = "  +1.18%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.89"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +1.71%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "55.14"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "2.549.86"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +2.18%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "73.16"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.27%  "
